# Applies the "Updated cryptos list" data refresh to Sheet1.
# Column D ("Price") values that look like plain decimal numbers must be
# forced to Text format first, otherwise Excel auto-converts them to
# numeric values (losing trailing zeros / exact text representation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.518.75'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '2.646.52'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.86'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.80'
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("E9").Value = '  +1.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.58'
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.369'
$ws.Range("E11").Value = '  +4.25%  '
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.47'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '3.125.39'
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").Value = '63.406.17'
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '2.647.51'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("E19").Value = '  +4.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.52'
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.92'
$ws.Range("E21").Value = '  +2.47%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.58'
$ws.Range("E23").Value = '  -3.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.67'
$ws.Range("E24").Value = '  -1.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.68'
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.11'
$ws.Range("E26").Value = '  +7.43%  '
$ws.Range("B27").Value = 'Bittensor'
$ws.Range("C27").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '567.72'
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("B28").Value = 'SuiNetwork'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.56'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.93'
$ws.Range("E31").Value = '  +1.34%  '
$ws.Range("E32").Value = '  +2.72%  '
$ws.Range("E33").Value = '  -3.81%  '
$ws.Range("D34").Value = '0.0₃0815'
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("E35").Value = '  +4.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.43'
$ws.Range("E36").Value = '  -4.00%  '
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("E39").Value = '  +4.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.09'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '168.71'
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.76'
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.12'
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0570'
$ws.Range("E45").Value = '  +2.58%  '
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("E47").Value = '  +3.19%  '
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.74'
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("E50").Value = '  +8.62%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.177'
$ws.Range("E51").Value = '  +2.32%  '
